# feat : working with image completed
# Append a new attachment record (Id=3) to the tracking sheet, mirroring
# the existing Id/Name/Location/FileSize/FileType columns in row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "C:\Users\Admin\Desktop\Best-Indian-Punjabi-Samosa-Recipe.jpg"
$ws.Range("C4").Value = "D:\work-place\flutter apps\projects\sharp\BisleriumCafeBackend\fyptempdocument\doc\1704735419442-5ce111d5-dd7c-41f3-b432-7abff9a14dd6..jpg"
$ws.Range("D4").Value = 0.070037841796875
$ws.Range("E4").Value = "IMAGE"
